$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended after the existing 3 data rows.
$ws.Range("A4").Value = "Hi"
# Force this numeric-looking value to be stored as text (matches the
# target workbook, where B4 is a shared string "1300003", not a number).
$ws.Range("B4").Value = "'1300003"
$ws.Range("C4").Value = "Bye"

# Touch the remaining columns of the row (no-op formatting) so the
# worksheet XML materializes explicit empty cells for D4:K4, matching
# the rest of the table's row layout (columns A:K).
$ws.Range("D4:K4").Font.Bold = $false
